$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: insert the newest price record (2021-09-16) as a new
# row right before the existing "Rabanito" block that starts at row 23,
# pushing all subsequent rows down by one (Excel's native Insert()
# semantics - matches the diff where rows 23-36 each now hold the data
# that used to live in the row right below them, and a new row 37 is
# appended with what used to be the last row's data).
$ws.Rows.Item(23).Insert()

$ws.Cells.Item(23, 1).Value = 10
$ws.Cells.Item(23, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(23, 3).Value = "La Araucanía"
$ws.Cells.Item(23, 4).Value = "2021-09-16"
$ws.Cells.Item(23, 5).Value = 9
$ws.Cells.Item(23, 6).Value = 300000001
$ws.Cells.Item(23, 7).Value = "Rabanito"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 20
$ws.Cells.Item(23, 11).Value = 7000
$ws.Cells.Item(23, 12).Value = 7000
$ws.Cells.Item(23, 13).Value = 7000
$ws.Cells.Item(23, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(23, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(23, 16).Value = 583
$ws.Cells.Item(23, 17).Value = 12
$ws.Cells.Item(23, 18).Value = "Hortaliza"
